# Sync attendance_reports: normalize the "Recorded By" (column G) ordering
# so the audit/system account name sorts after the human/service email
# address it was logged alongside.
#
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System, system" -> "backup@backdoor.com, system, System"
#
# Any other combination (e.g. a lone "System", "dnasr281@gmail.com" by
# itself, "admin@admin.com, System", "backup@backdoor.com, System") is left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$col = 7  # column G = "Recorded By"
$changed = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
        $changed++
    } elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "backup@backdoor.com, system, System"
        $changed++
    }
}

Write-Output "Updated $changed cell(s) in column G"
